$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.088.13'
$ws.Range("E2").Value = '  -1.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.306.11'
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.91'
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.16'
$ws.Range("E6").Value = '  -1.55%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.611'
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.84'
$ws.Range("E10").Value = '  -4.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0905'
$ws.Range("E11").Value = '  -2.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.48'
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.979'
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.51'
$ws.Range("E15").Value = '  -3.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.654.57'
$ws.Range("E16").Value = '  -1.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.317.72'
$ws.Range("E17").Value = '  -1.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.061.58'
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.75'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.87'
$ws.Range("E21").Value = '  -4.27%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '283.04'
$ws.Range("E22").Value = '  +9.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.57'
$ws.Range("E23").Value = '  -1.33%  '
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.02'
$ws.Range("E25").Value = '  +6.16%  '
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.99'
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.98'
$ws.Range("E28").Value = '  -3.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '23.39'
$ws.Range("E29").Value = '  +1.61%  '
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '165.69'
$ws.Range("E31").Value = '  -5.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '35.61'
$ws.Range("E32").Value = '  -2.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0885'
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.92'
$ws.Range("E35").Value = '  -3.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.133'
$ws.Range("E36").Value = '  +1.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.68'
$ws.Range("E38").Value = '  +1.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.92'
$ws.Range("E39").Value = '  +8.24%  '
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.65'
$ws.Range("E41").Value = '  -3.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '102.60'
$ws.Range("E42").Value = '  +20.69%  '
$ws.Range("E43").Value = '  +1.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.86'
$ws.Range("E44").Value = '  -1.49%  '
$ws.Range("E45").Value = '  -4.72%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '117.48'
$ws.Range("E47").Value = '  +2.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.16'
$ws.Range("E48").Value = '  +1.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '78.87'
$ws.Range("E49").Value = '  +6.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.15'
$ws.Range("E50").Value = '  +0.32%  '
$ws.Range("E51").Value = '  -2.68%  '
